$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.418174028396606
$ws.Range("B1").Value = 6.760646820068359
$ws.Range("C1").Value = 5.919130802154541
$ws.Range("D1").Value = 2.275228977203369
$ws.Range("E1").Value = 1.431994199752808
